# The search results table gained one more appended row (row 6), recording
# another "Random" method run whose results the user chose to keep.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 6

# Copy the number format (date style) from the row above so the new date
# cell renders the same way as the existing ones.
$ws.Range("A5").Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A" + $newRow).Value = 42602.584120370368
$ws.Range("B" + $newRow).Value = "Random"
$ws.Range("C" + $newRow).Value = 0
$ws.Range("D" + $newRow).Value = 0
$ws.Range("E" + $newRow).Value = 0
$ws.Range("F" + $newRow).Value = 0
$ws.Range("G" + $newRow).Value = 0
$ws.Range("H" + $newRow).Value = 85
$ws.Range("I" + $newRow).Value = 15
$ws.Range("J" + $newRow).Value = 0
$ws.Range("K" + $newRow).Value = 0
$ws.Range("L" + $newRow).Value = 86
$ws.Range("M" + $newRow).Value = 14
